$d = $word.ActiveDocument

# Sanity-check: make sure the anchor paragraph (last paragraph of the
# document body) still contains the expected trailing sentence before we
# append new content after it.
$checkRng = $d.Content
$found = $checkRng.Find.Execute(
    "do Branch corrente pode ser chamado pela palavra HEAD. ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Anchor paragraph text not found; aborting edit."
}

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Two new empty paragraphs, matching the justification/size formatting of
# the surrounding text but with no run content at all.
$emptyParaXml = '<w:p ' + $wNs + '>' +
    '<w:pPr>' +
        '<w:jc w:val="both"/>' +
        '<w:rPr>' +
            '<w:sz w:val="28"/>' +
            '<w:szCs w:val="28"/>' +
        '</w:rPr>' +
    '</w:pPr>' +
'</w:p>'

# Final new paragraph: underlined paragraph mark, with a normal (non
# underlined) run carrying the new sentence.
$noteText = "Essa trata-se de uma modificação não necessária. "
$noteParaXml = '<w:p ' + $wNs + '>' +
    '<w:pPr>' +
        '<w:jc w:val="both"/>' +
        '<w:rPr>' +
            '<w:sz w:val="28"/>' +
            '<w:szCs w:val="28"/>' +
            '<w:u w:val="single"/>' +
        '</w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
        '<w:rPr>' +
            '<w:sz w:val="28"/>' +
            '<w:szCs w:val="28"/>' +
        '</w:rPr>' +
        '<w:t xml:space="preserve">' + $noteText + '</w:t>' +
    '</w:r>' +
'</w:p>'

# Insert all three new paragraphs in one go right at the very end of the
# document (immediately before the final section properties), so they
# land directly after the "...palavra HEAD. " paragraph.
$combinedXml = $emptyParaXml + $emptyParaXml + $noteParaXml

$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertXML($combinedXml)
